$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(4, 1).Value = "'2025-10-13"
$ws.Cells.Item(4, 2).Value = "'12180.00"
